# Hot fix edited tab names to test new generic tab implementation
$wb = $excel.ActiveWorkbook

# Rename the "Project" group's child tabs to be prefixed with "Project - "
$wb.Worksheets.Item("Contact").Name = "Project - Contact"
$wb.Worksheets.Item("Publications").Name = "Project - Publications"
$wb.Worksheets.Item("Funder").Name = "Project - Funder"

# Move the active/selected tab from "Cell suspension" back to "Project"
$wb.Worksheets.Item("Project").Activate()
